$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = 45993
$ws.Cells.Item(2, 2).Value = 11022.8167045037
$ws.Cells.Item(2, 3).Value = 10175.3454363828
$ws.Cells.Item(2, 4).Value = 17064.26
$ws.Cells.Item(2, 5).Value = 6974.75794216093
$ws.Cells.Item(2, 6).Value = 3.57680743932254

$ws.Cells.Item(3, 1).Value = 45994
$ws.Cells.Item(3, 2).Value = 11101.5213286104
$ws.Cells.Item(3, 3).Value = 10236.0942462511
$ws.Cells.Item(3, 4).Value = 11232.26
$ws.Cells.Item(3, 5).Value = 7039.07391346623
$ws.Cells.Item(3, 6).Value = 251.78783998822

$ws.Cells.Item(4, 1).Value = 45995
$ws.Cells.Item(4, 2).Value = 10948.2721276786
$ws.Cells.Item(4, 3).Value = 10126.0612947905
$ws.Cells.Item(4, 4).Value = 11232.26
$ws.Cells.Item(4, 5).Value = 6928.86515349851
$ws.Cells.Item(4, 6).Value = 242.611102012043

$ws.Cells.Item(5, 1).Value = 45996
$ws.Cells.Item(5, 2).Value = 10750.3021268293
$ws.Cells.Item(5, 3).Value = 9371.69155106782
$ws.Cells.Item(5, 4).Value = 11232.26
$ws.Cells.Item(5, 5).Value = 6781.48970570058
$ws.Cells.Item(5, 6).Value = 205.038385698683

$ws.Cells.Item(6, 1).Value = 45997
$ws.Cells.Item(6, 2).Value = 8388.99607296656
$ws.Cells.Item(6, 3).Value = 8898.51798427265
$ws.Cells.Item(6, 4).Value = 11232.26
$ws.Cells.Item(6, 5).Value = 7096.26940993017
$ws.Cells.Item(6, 6).Value = 198.438641425118

$ws.Cells.Item(7, 1).Value = 45998
$ws.Cells.Item(7, 2).Value = 8289.54708377905
$ws.Cells.Item(7, 3).Value = 8696.21068854051
$ws.Cells.Item(7, 4).Value = 11232.26
$ws.Cells.Item(7, 5).Value = 7081.20396894763
$ws.Cells.Item(7, 6).Value = 189.381444062005

$ws.Cells.Item(8, 1).Value = 45999
$ws.Cells.Item(8, 2).Value = 8449.17139773619
$ws.Cells.Item(8, 3).Value = 8760.74675699143
$ws.Cells.Item(8, 4).Value = 11232.26
$ws.Cells.Item(8, 5).Value = 7514.74508513084
$ws.Cells.Item(8, 6).Value = 210.134660088428

$ws.Cells.Item(9, 1).Value = 46000
$ws.Cells.Item(9, 2).Value = 9791.8538926876
$ws.Cells.Item(9, 3).Value = 9773.12445858814
$ws.Cells.Item(9, 4).Value = 11232.26
$ws.Cells.Item(9, 5).Value = 7938.45608770007
$ws.Cells.Item(9, 6).Value = 269.971689428675

$ws.Cells.Item(10, 1).Value = 46001
$ws.Cells.Item(10, 2).Value = 9791.8538926876
$ws.Cells.Item(10, 3).Value = 9449.4031374742
$ws.Cells.Item(10, 4).Value = 11232.26
$ws.Cells.Item(10, 5).Value = 7938.45608770007
$ws.Cells.Item(10, 6).Value = 256.483301048928

$ws.Cells.Item(11, 1).Value = 46002
$ws.Cells.Item(11, 2).Value = 9791.8538926876
$ws.Cells.Item(11, 3).Value = 9444.10935814389
$ws.Cells.Item(11, 4).Value = 11232.26
$ws.Cells.Item(11, 5).Value = 7938.45608770007
$ws.Cells.Item(11, 6).Value = 256.262726910165

$ws.Cells.Item(12, 1).Value = 46003
$ws.Cells.Item(12, 2).Value = 9791.8538926876
$ws.Cells.Item(12, 3).Value = 8702.67441728799
$ws.Cells.Item(12, 4).Value = 11232.26
$ws.Cells.Item(12, 5).Value = 7938.45608770007
$ws.Cells.Item(12, 6).Value = 225.369604374502

$ws.Cells.Item(13, 1).Value = 46004
$ws.Cells.Item(13, 2).Value = 8551.33924250961
$ws.Cells.Item(13, 3).Value = 8479.87149811022
$ws.Cells.Item(13, 4).Value = 11232.26
$ws.Cells.Item(13, 5).Value = 7530.17038790506
$ws.Cells.Item(13, 6).Value = 199.074245250637

$ws.Cells.Item(14, 1).Value = 46005
$ws.Cells.Item(14, 2).Value = 8449.17139773619
$ws.Cells.Item(14, 3).Value = 8423.65194407213
$ws.Cells.Item(14, 4).Value = 11232.26
$ws.Cells.Item(14, 5).Value = 7514.82167260086
$ws.Cells.Item(14, 6).Value = 196.092234028041

$ws.Cells.Item(15, 1).Value = 46006
$ws.Cells.Item(15, 2).Value = 9832.03344432964
$ws.Cells.Item(15, 3).Value = 9157.72779691708
$ws.Cells.Item(15, 4).Value = 11232.26
$ws.Cells.Item(15, 5).Value = 8217.9042789547
$ws.Cells.Item(15, 6).Value = 255.973836494657

